# "Stop donation on already reciveved requests"
#
# The canonical-OOXML diff for this commit shows slide 5 ("Key Points" -
# the crowdfunding/donation summary slide) being removed from the deck.
# Every other change in the diff (notesMaster/handoutMaster rId shift,
# slide5.xml/slide6.xml/slide7.xml file renumbering, the notesSlide5/6
# merge, the cached slide-number field text) is the mechanical fallout of
# that single deletion once the package is renumbered/exported - so the
# one deliberate edit to replay here is removing that slide.

$p = $ppt.ActivePresentation

# Find the "Key Points" slide defensively (it is slide 5 in the original
# deck) instead of hard-coding an index, in case slides have already
# shifted.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -like "Key Points*") {
                $targetIndex = $i
                break
            }
        }
    }
    if ($targetIndex -ne -1) { break }
}

if ($targetIndex -eq -1) {
    # Fallback: original known position.
    $targetIndex = 5
}

$p.Slides.Item($targetIndex).Delete()

# Best-effort: the diff also shows the notesMaster/handoutMaster
# "datetimeFigureOut" cached field text moving from 4/2/2022 to
# 6/11/2022 (an auto date field re-cached on a later save). Try to
# refresh it on both masters; harmless if unsupported.
try {
    $nm = $p.NotesMaster
    for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
        $sh = $nm.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "6/11/2022"
        }
    }
} catch {
}

try {
    $hm = $p.HandoutMaster
    for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
        $sh = $hm.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "6/11/2022"
        }
    }
} catch {
}
